# Fix the age-range header typo "5a 14" -> "5 a 14" (missing space)
# across all the virus-specific sheets (cells F1 and N1).

$wb = $excel.ActiveWorkbook

$sheetNames = @("VRS", "Ad", "Parainfluenza", "Inf A", "Inf B", "Metapnemovirus")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F1").Value = "5 a 14"
    $ws.Range("N1").Value = "5 a 14"
}
